# Generate Report for Handoff
# Rename the source file's generated UUID-based filename/hash from
#   c46f4aee-e486-42fc-977d-ba2d660ddee2 / c15269bb08f6ab0a3aeb8d6b4c1c6adaf164d513
# to
#   49e75141-5af2-4885-aca5-2eac0a7b72b5 / c4f2c270380f0c69004b8fd0813c2898526f6860
# across the Overview / zh-cn / de-de sheets, and bump the handoff timestamps.

$wb = $excel.ActiveWorkbook

$newGuid = "49e75141-5af2-4885-aca5-2eac0a7b72b5"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet ---
$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("B2").Value = "e2e\$newGuid.md"
$wsOverview.Range("G2").Value = "2016-09-05 05:04:00"

# --- zh-cn sheet ---
$wsZhCn.Range("A2").Value = "$newGuid.md"
$wsZhCn.Range("G2").Value = "$newGuid.c4f2c270380f0c69004b8fd0813c2898526f6860.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-09-05 05:03:55"

# --- de-de sheet ---
$wsDeDe.Range("A2").Value = "$newGuid.md"
$wsDeDe.Range("G2").Value = "$newGuid.c4f2c270380f0c69004b8fd0813c2898526f6860.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-09-05 05:04:00"
